# Update simulation results for the 380 kV case (pl_mw.xlsx, Case_3_107).
# Columns B,C,D,F,G,I,K,L,M,N,O for rows 2-25 get new computed values;
# columns A,E,H,J (and the header row) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4693012460242301
$ws.Range("C2").Value = 0.1639138001180811
$ws.Range("D2").Value = 0.04385257605817827
$ws.Range("F2").Value = 0.8891880782175861
$ws.Range("G2").Value = 0.002455405213520915
$ws.Range("I2").Value = 0.8624332825623817
$ws.Range("K2").Value = 0.2716488505277255
$ws.Range("L2").Value = 0.2971665877106915
$ws.Range("M2").Value = 0.173335037061868
$ws.Range("N2").Value = 1.835328718841951
$ws.Range("O2").Value = 3.154250357694863

$ws.Range("B3").Value = 0.4344328004809768
$ws.Range("C3").Value = 0.1626536677307229
$ws.Range("D3").Value = 0.04152777417815656
$ws.Range("F3").Value = 0.8893646934460406
$ws.Range("G3").Value = 0.002457605324409473
$ws.Range("I3").Value = 0.8680153509411674
$ws.Range("K3").Value = 0.2392153479422348
$ws.Range("L3").Value = 0.2937218382855562
$ws.Range("M3").Value = 0.1660851520942295
$ws.Range("N3").Value = 1.85281748240822
$ws.Range("O3").Value = 3.16748342244459

$ws.Range("B4").Value = 0.4131539817925614
$ws.Range("C4").Value = 0.1618735999780583
$ws.Range("D4").Value = 0.04008455143407019
$ws.Range("F4").Value = 0.8899414794722986
$ws.Range("G4").Value = 0.002459029671682291
$ws.Range("I4").Value = 0.8718179859160138
$ws.Range("K4").Value = 0.2192884055372843
$ws.Range("L4").Value = 0.2917552650435411
$ws.Range("M4").Value = 0.1617081404684697
$ws.Range("N4").Value = 1.864105693901423
$ws.Range("O4").Value = 3.177128910114277

$ws.Range("B5").Value = 0.4045161232638748
$ws.Range("C5").Value = 0.1615541415053414
$ws.Range("D5").Value = 0.03949248099451808
$ws.Range("F5").Value = 0.8902944053646422
$ws.Range("G5").Value = 0.002459628632768455
$ws.Range("I5").Value = 0.8734620111411644
$ws.Range("K5").Value = 0.2111652990645183
$ws.Range("L5").Value = 0.2909913055576325
$ws.Range("M5").Value = 0.1599433255534031
$ws.Range("N5").Value = 1.868844079111827
$ws.Range("O5").Value = 3.181442092806961

$ws.Range("B6").Value = 0.4030838505111092
$ws.Range("C6").Value = 0.1615010011640337
$ws.Range("D6").Value = 0.03939393055211582
$ws.Range("F6").Value = 0.8903601314531215
$ws.Range("G6").Value = 0.00245972921024017
$ws.Range("I6").Value = 0.8737407050485615
$ws.Range("K6").Value = 0.2098163130261099
$ws.Range("L6").Value = 0.2908667142784012
$ws.Range("M6").Value = 0.159651421945334
$ws.Range("N6").Value = 1.869639241280423
$ws.Range("O6").Value = 3.182181407585517

$ws.Range("B7").Value = 0.4130373524406536
$ws.Range("C7").Value = 0.1618692979965672
$ws.Range("D7").Value = 0.04007658250915114
$ws.Range("F7").Value = 0.889945761710436
$ws.Range("G7").Value = 0.002459037674144575
$ws.Range("I7").Value = 0.8718397754196232
$ws.Range("K7").Value = 0.2191788648651567
$ws.Range("L7").Value = 0.2917448103267901
$ws.Range("M7").Value = 0.1616842630470217
$ws.Range("N7").Value = 1.864169037210865
$ws.Range("O7").Value = 3.177185529845701

$ws.Range("B8").Value = 0.4572518628789339
$ws.Range("C8").Value = 0.1634806342935775
$ws.Range("D8").Value = 0.04305427423658159
$ws.Range("F8").Value = 0.8891518354489989
$ws.Range("G8").Value = 0.002456148595264202
$ws.Range("I8").Value = 0.8642801427003661
$ws.Range("K8").Value = 0.260468670475916
$ws.Range("L8").Value = 0.2959480652904745
$ws.Range("M8").Value = 0.1708199022267678
$ws.Range("N8").Value = 1.841244706617411
$ws.Range("O8").Value = 3.158497730356359

$ws.Range("B9").Value = 0.5449700801460722
$ws.Range("C9").Value = 0.1665894004648649
$ws.Range("D9").Value = 0.04876746394604226
$ws.Range("F9").Value = 0.8913071903246319
$ws.Range("G9").Value = 0.00245106363058039
$ws.Range("I9").Value = 0.8524301832512755
$ws.Range("K9").Value = 0.3413210640961495
$ws.Range("L9").Value = 0.3053658773178682
$ws.Range("M9").Value = 0.1893207779654063
$ws.Range("N9").Value = 1.800652218175659
$ws.Range("O9").Value = 3.133904399788292

$ws.Range("B10").Value = 0.6100114461708017
$ws.Range("C10").Value = 0.1688414877304822
$ws.Range("D10").Value = 0.05288736332612132
$ws.Range("F10").Value = 0.8951502148792443
$ws.Range("G10").Value = 0.002447678183321219
$ws.Range("I10").Value = 0.8455339334288468
$ws.Range("K10").Value = 0.4006356261407404
$ws.Range("L10").Value = 0.3129984473619629
$ws.Range("M10").Value = 0.2032656492345026
$ws.Range("N10").Value = 1.77348404253109
$ws.Range("O10").Value = 3.123173600501474

$ws.Range("B11").Value = 0.6397250623340085
$ws.Range("C11").Value = 0.1698589322877098
$ws.Range("D11").Value = 0.0547446285953086
$ws.Range("F11").Value = 0.8973884691440475
$ws.Range("G11").Value = 0.002446213431851733
$ws.Range("I11").Value = 0.8427889706567342
$ws.Range("K11").Value = 0.4275971392585518
$ws.Range("L11").Value = 0.3166249832002848
$ws.Range("M11").Value = 0.2096850714121103
$ws.Range("N11").Value = 1.7617001956509
$ws.Range("O11").Value = 3.119883315057251

$ws.Range("B12").Value = 0.6509944169608559
$ws.Range("C12").Value = 0.1702431817431886
$ws.Range("D12").Value = 0.05544547633562047
$ws.Range("F12").Value = 0.8983064190422283
$ws.Range("G12").Value = 0.002445669542626242
$ws.Range("I12").Value = 0.8418058614552706
$ws.Range("K12").Value = 0.4378033450510941
$ws.Range("L12").Value = 0.3180203912593385
$ws.Range("M12").Value = 0.2121267286354325
$ws.Range("N12").Value = 1.757320626269742
$ws.Range("O12").Value = 3.118865978662029

$ws.Range("B13").Value = 0.6485665941215473
$ws.Range("C13").Value = 0.1701604731196795
$ws.Range("D13").Value = 0.05529464602363987
$ws.Range("F13").Value = 0.8981055938507083
$ws.Range("G13").Value = 0.002445786200327621
$ws.Range("I13").Value = 0.8420150862201936
$ws.Range("K13").Value = 0.4356054214470078
$ws.Range("L13").Value = 0.3177188829305209
$ws.Range("M13").Value = 0.2116003977308765
$ws.Range("N13").Value = 1.758260165660428
$ws.Range("O13").Value = 3.119074915219443

$ws.Range("B14").Value = 0.6406518529285847
$ws.Range("C14").Value = 0.1698905655840974
$ws.Range("D14").Value = 0.05480233719227812
$ws.Range("F14").Value = 0.897462579636418
$ws.Range("G14").Value = 0.002446168470010052
$ws.Range("I14").Value = 0.8427069606227775
$ws.Range("K14").Value = 0.4284368837466559
$ws.Range("L14").Value = 0.3167393415127862
$ws.Range("M14").Value = 0.2098857331676385
$ws.Range("N14").Value = 1.761338227933139
$ws.Range("O14").Value = 3.119795037323513

$ws.Range("B15").Value = 0.6358060986674161
$ws.Range("C15").Value = 0.1697251040178713
$ws.Range("D15").Value = 0.05450046282732757
$ws.Range("F15").Value = 0.8970778764565281
$ws.Range("G15").Value = 0.002446404023759196
$ws.Range("I15").Value = 0.8431380905684662
$ws.Range("K15").Value = 0.4240454718065223
$ws.Range("L15").Value = 0.3161422217764454
$ws.Range("M15").Value = 0.2088368488115293
$ws.Range("N15").Value = 1.763234402873412
$ws.Range("O15").Value = 3.12026590073512

$ws.Range("B16").Value = 0.6080720698229811
$ws.Range("C16").Value = 0.1687748517492622
$ws.Range("D16").Value = 0.05276564421968288
$ws.Range("F16").Value = 0.8950137924311861
$ws.Range("G16").Value = 0.002447775419458871
$ws.Range("I16").Value = 0.8457212108673033
$ws.Range("K16").Value = 0.398873161497761
$ws.Range("L16").Value = 0.3127645447171545
$ws.Range("M16").Value = 0.2028476385868387
$ws.Range("N16").Value = 1.774265718632814
$ws.Range("O16").Value = 3.123420631309784

$ws.Range("B17").Value = 0.5910899355014294
$ws.Range("C17").Value = 0.1681900835227452
$ws.Range("D17").Value = 0.05169704110734585
$ws.Range("F17").Value = 0.8938729874889972
$ws.Range("G17").Value = 0.002448635979121038
$ws.Range("I17").Value = 0.8474062845533865
$ws.Range("K17").Value = 0.3834250234456817
$ws.Range("L17").Value = 0.3107319409082692
$ws.Range("M17").Value = 0.1991927724717044
$ws.Range("N17").Value = 1.781180382078848
$ws.Range("O17").Value = 3.125763367752114

$ws.Range("B18").Value = 0.5813341577324707
$ws.Range("C18").Value = 0.1678530788762274
$ws.Range("D18").Value = 0.05108082008570847
$ws.Range("F18").Value = 0.8932629576303768
$ws.Range("G18").Value = 0.00244913804156863
$ws.Range("I18").Value = 0.8484124105488924
$ws.Range("K18").Value = 0.3745377193916397
$ws.Range("L18").Value = 0.3095773831793593
$ws.Range("M18").Value = 0.1970977362239665
$ws.Range("N18").Value = 1.785211642679072
$ws.Range("O18").Value = 3.127260631348662

$ws.Range("B19").Value = 0.5780330824298971
$ws.Range("C19").Value = 0.1677388619674431
$ws.Range("D19").Value = 0.05087190622734283
$ws.Range("F19").Value = 0.8930643379533123
$ws.Range("G19").Value = 0.002449309250288654
$ws.Range("I19").Value = 0.8487594093099027
$ws.Range("K19").Value = 0.3715283128749434
$ws.Range("L19").Value = 0.3091889701711068
$ws.Range("M19").Value = 0.1963896246426984
$ws.Range("N19").Value = 1.786585854400707
$ws.Range("O19").Value = 3.127793310345936

$ws.Range("B20").Value = 0.5928964862144142
$ws.Range("C20").Value = 0.1682524017162237
$ws.Range("D20").Value = 0.05181096037130573
$ws.Range("F20").Value = 0.8939896544329571
$ws.Range("G20").Value = 0.002448543637677131
$ws.Range("I20").Value = 0.8472230852283111
$ws.Range("K20").Value = 0.3850697090171593
$ws.Range("L20").Value = 0.3109468106456745
$ws.Range("M20").Value = 0.199581100651713
$ws.Range("N20").Value = 1.78043870205565
$ws.Range("O20").Value = 3.125498478798846

$ws.Range("B21").Value = 0.6429761358796497
$ws.Range("C21").Value = 0.1699698722308085
$ws.Range("D21").Value = 0.05494700717557777
$ws.Range("F21").Value = 0.8976495394030124
$ws.Range("G21").Value = 0.002446055896028597
$ws.Range("I21").Value = 0.8425022114184841
$ws.Range("K21").Value = 0.4305425570783257
$ws.Range("L21").Value = 0.3170264570314743
$ws.Range("M21").Value = 0.2103890807946271
$ws.Range("N21").Value = 1.760431881324301
$ws.Range("O21").Value = 3.119577316758125

$ws.Range("B22").Value = 0.6758075040796143
$ws.Range("C22").Value = 0.1710862971118274
$ws.Range("D22").Value = 0.05698225284162106
$ws.Range("F22").Value = 0.900451620074314
$ws.Range("G22").Value = 0.002444492826886979
$ws.Range("I22").Value = 0.8397452691910914
$ws.Range("K22").Value = 0.4602408519615153
$ws.Range("L22").Value = 0.3211287334556516
$ws.Range("M22").Value = 0.2175153794806945
$ws.Range("N22").Value = 1.74783842039111
$ws.Range("O22").Value = 3.117039997542349

$ws.Range("B23").Value = 0.6582757236768941
$ws.Range("C23").Value = 0.1704910002480844
$ws.Range("D23").Value = 0.05589732567531769
$ws.Range("F23").Value = 0.8989186015224533
$ws.Range("G23").Value = 0.002445321334647209
$ws.Range("I23").Value = 0.841186665648344
$ws.Range("K23").Value = 0.4443923910814078
$ws.Range("L23").Value = 0.3189275092922941
$ws.Range("M23").Value = 0.2137062537865972
$ws.Range("N23").Value = 1.754515663230136
$ws.Range("O23").Value = 3.118272351694117

$ws.Range("B24").Value = 0.5920797212958746
$ws.Range("C24").Value = 0.1682242301939141
$ws.Range("D24").Value = 0.0517594632852294
$ws.Range("F24").Value = 0.8939367665310129
$ws.Range("G24").Value = 0.00244858536254061
$ws.Range("I24").Value = 0.8473057932824481
$ws.Range("K24").Value = 0.3843261652077672
$ws.Range("L24").Value = 0.3108496243814614
$ws.Range("M24").Value = 0.1994055181744585
$ws.Range("N24").Value = 1.780773841380611
$ws.Range("O24").Value = 3.125617766615477

$ws.Range("B25").Value = 0.5211338494018776
$ws.Range("C25").Value = 0.1657539455850454
$ws.Range("D25").Value = 0.04723546525708144
$ws.Range("F25").Value = 0.8903270580901861
$ws.Range("G25").Value = 0.002452377455620935
$ws.Range("I25").Value = 0.8553178266520476
$ws.Range("K25").Value = 0.3194625550210901
$ws.Range("L25").Value = 0.3026926166362927
$ws.Range("M25").Value = 0.1842535812465869
$ws.Range("N25").Value = 1.811167007235759
$ws.Range("O25").Value = 3.139268197544084
